$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hlookup")
$ws.Activate()

# Row 7 (dmart): HLOOKUP against $D$15:$J$17, row 2 of the table
$ws.Range("E7").Formula = '=HLOOKUP(E6,$D$15:$J$17,2,FALSE)'
$ws.Range("F7:J7").Formula = '=HLOOKUP(F6,$D$15:$J$17,2,FALSE)'

# Row 10 (swiggy instamart): HLOOKUP against $D$15:$J$17, row 3 of the table
$ws.Range("E10").Formula = '=HLOOKUP(E6,$D$15:$J$17,3,FALSE)'
$ws.Range("F10:J10").Formula = '=HLOOKUP(F6,$D$15:$J$17,3,FALSE)'

# Update the selection to match the edited range
$ws.Range("E10:J10").Select()
